$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item(1))
$names = $ws.PageSetup | Get-Member | Out-String
Write-Host $names
